$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to be treated as plain text so that
# values like "1.00", "65.745.89" or "0.0854" keep their exact
# original formatting instead of being auto-coerced to numbers.
foreach ($row in @(2,3,4,5,6,7,9,10,13,15,17,18,19,20,21,22,23,24,26,28,29,31,34,35,39,40,41,42,43,44,45,46,47,49,50,51)) {
    $ws.Cells.Item($row, 4).NumberFormat = "@"
}

$ws.Range("D2").Value = '65.745.89'
$ws.Range("E2").Value = '  -1.50%  '
$ws.Range("D3").Value = '3.416.38'
$ws.Range("E3").Value = '  -1.91%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = '596.82'
$ws.Range("E5").Value = '  -1.03%  '
$ws.Range("D6").Value = '142.51'
$ws.Range("E6").Value = '  -3.51%  '
$ws.Range("D7").Value = '3.415.07'
$ws.Range("E7").Value = '  -1.85%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").Value = '0.471'
$ws.Range("E9").Value = '  -2.44%  '
$ws.Range("D10").Value = '7.99'
$ws.Range("E10").Value = '  +5.63%  '
$ws.Range("E11").Value = '  -5.43%  '
$ws.Range("E12").Value = '  -4.27%  '
$ws.Range("D13").Value = '3.991.07'
$ws.Range("E13").Value = '  -2.00%  '
$ws.Range("E14").Value = '  -6.20%  '
$ws.Range("D15").Value = '29.62'
$ws.Range("E15").Value = '  -5.80%  '
$ws.Range("E16").Value = '  -0.70%  '
$ws.Range("D17").Value = '3.412.07'
$ws.Range("E17").Value = '  -1.95%  '
$ws.Range("D18").Value = '65.672.89'
$ws.Range("E18").Value = '  -1.69%  '
$ws.Range("D19").Value = '10.39'
$ws.Range("E19").Value = '  +3.29%  '
$ws.Range("D20").Value = '6.12'
$ws.Range("E20").Value = '  -4.96%  '
$ws.Range("D21").Value = '14.59'
$ws.Range("E21").Value = '  -4.96%  '
$ws.Range("D22").Value = '415.29'
$ws.Range("E22").Value = '  -5.24%  '
$ws.Range("D23").Value = '0.578'
$ws.Range("E23").Value = '  -5.05%  '
$ws.Range("D24").Value = '77.25'
$ws.Range("E24").Value = '  -2.76%  '
$ws.Range("E25").Value = '  +0.08%  '
$ws.Range("D26").Value = '3.548.57'
$ws.Range("E26").Value = '  -2.00%  '
$ws.Range("E27").Value = '  -9.04%  '
$ws.Range("D28").Value = '9.25'
$ws.Range("E28").Value = '  -5.72%  '
$ws.Range("D29").Value = '7.83'
$ws.Range("E29").Value = '  -6.68%  '
$ws.Range("E30").Value = '  -2.62%  '
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  -0.13%  '
$ws.Range("E32").Value = '  -4.75%  '
$ws.Range("E33").Value = '  -8.16%  '
$ws.Range("D34").Value = '24.55'
$ws.Range("E34").Value = '  -3.23%  '
$ws.Range("D35").Value = '3.409.64'
$ws.Range("E35").Value = '  -1.86%  '
$ws.Range("E37").Value = '  -6.35%  '
$ws.Range("E38").Value = '  -8.84%  '
$ws.Range("D39").Value = '7.52'
$ws.Range("E39").Value = '  -5.08%  '
$ws.Range("D40").Value = '0.998'
$ws.Range("E40").Value = '  -0.13%  '
$ws.Range("D41").Value = '168.52'
$ws.Range("E41").Value = '  -4.85%  '
$ws.Range("D42").Value = '0.0854'
$ws.Range("E42").Value = '  -3.39%  '
$ws.Range("D43").Value = '0.874'
$ws.Range("E43").Value = '  -2.09%  '
$ws.Range("D44").Value = '5.04'
$ws.Range("E44").Value = '  -7.16%  '
$ws.Range("D45").Value = '1.90'
$ws.Range("E45").Value = '  -10.71%  '
$ws.Range("D46").Value = '45.39'
$ws.Range("E46").Value = '  -2.08%  '
$ws.Range("D47").Value = '26.20'
$ws.Range("E47").Value = '  -9.13%  '
$ws.Range("E48").Value = '  -4.18%  '
$ws.Range("D49").Value = '7.05'
$ws.Range("E49").Value = '  -5.46%  '
$ws.Range("D50").Value = '2.26'
$ws.Range("E50").Value = '  -6.74%  '
$ws.Range("D51").Value = '0.919'
$ws.Range("E51").Value = '  -6.22%  '

# Restore the default cell style on the Price cells we reformatted as text,
# so no stray number-format/style residue is left behind.
foreach ($row in @(2,3,4,5,6,7,9,10,13,15,17,18,19,20,21,22,23,24,26,28,29,31,34,35,39,40,41,42,43,44,45,46,47,49,50,51)) {
    $ws.Cells.Item($row, 4).Style = "Normal"
}
